$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refreshed "last updated" timestamp (sharedStrings text) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 06:52"

# --- Reorder countries: Tailandia's case count overtook Argelia's, so it
#     now sorts one row higher. Swap the two country-name labels and give
#     each row the refreshed data for its (possibly new) country. ---
$ws.Range("A58").Value = "Tailandia"
$ws.Range("B58").Value = 2826
$ws.Range("C58").Value = 15
$ws.Range("D58").Value = 2352
$ws.Range("E58").Value = 425
$ws.Range("F58").Value = 61
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 49

$ws.Range("A59").Value = "Argelia"
$ws.Range("B59").Value = 2811
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 1152
$ws.Range("E59").Value = 1267
$ws.Range("F59").Value = 40
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 392

# --- Reorder countries: Kirguistan's case count overtook Bolivia,
#     Albania and Burkina Faso, so it jumps to the top of that block;
#     the other three each shift down one row keeping their own data. ---
$ws.Range("A98").Value = "Kirguistan"
$ws.Range("B98").Value = 612
$ws.Range("C98").Value = 22
$ws.Range("D98").Value = 254
$ws.Range("E98").Value = 351
$ws.Range("F98").Value = 5
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 7

$ws.Range("A99").Value = "Bolivia"
$ws.Range("B99").Value = 609
$ws.Range("C99").Value = 11
$ws.Range("D99").Value = 44
$ws.Range("E99").Value = 528
$ws.Range("F99").Value = 3
$ws.Range("G99").Value = 3
$ws.Range("H99").Value = 37

$ws.Range("A100").Value = "Albania"
$ws.Range("B100").Value = 609
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 345
$ws.Range("E100").Value = 238
$ws.Range("F100").Value = 10
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 26

$ws.Range("A101").Value = "Burkina Faso"
$ws.Range("B101").Value = 600
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 362
$ws.Range("E101").Value = 200
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 38

# --- Plain data refreshes (no reordering) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 819175
$ws.Range("C4").Value = 431
$ws.Range("E4").Value = 690859
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 45343

# Row 20: India
$ws.Range("B20").Value = 20111
$ws.Range("C20").Value = 31
$ws.Range("E20").Value = 15491

# Row 32: Chile
$ws.Range("B32").Value = 9749
$ws.Range("C32").Value = 184
$ws.Range("D32").Value = 2156
$ws.Range("E32").Value = 7384
$ws.Range("F32").Value = 58
$ws.Range("G32").Value = 8
$ws.Range("H32").Value = 209

# Row 43: Indonesia
$ws.Range("D43").Value = 4920
$ws.Range("E43").Value = 1653

# Row 62: Hungria
$ws.Range("B62").Value = 2168
$ws.Range("C62").Value = 70
$ws.Range("D62").Value = 295
$ws.Range("E62").Value = 1660
